# Update cryptos list - GitHub Actions scheduled refresh
# (Sat Jul 29 15:35:51 UTC 2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not auto-coerced to a number),
# while leaving the cell's style/number-format exactly as it was before.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.299.28"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.873.25"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - XRP
Set-TextValue "D5" "0.7083"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6 - BNB
Set-TextValue "D6" "241.81"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - USDC
Set-TextValue "D7" "1.000"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.07804"
$ws.Range("E8").Value = "  +0.98%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.3107"
$ws.Range("E9").Value = "  -0.30%  "

# Row 10 - Solana
Set-TextValue "D10" "25.06"
$ws.Range("E10").Value = "  -1.48%  "

# Row 11 - TRON
Set-TextValue "D11" "0.08398"
$ws.Range("E11").Value = "  +0.08%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.871.40"
$ws.Range("E12").Value = "  -1.42%  "

# Row 13 - Polkadot
Set-TextValue "D13" "5.241"
$ws.Range("E13").Value = "  -0.40%  "

# Row 14 - Polygon
Set-TextValue "D14" "0.7175"
$ws.Range("E14").Value = "  -0.19%  "

# Row 15 - Litecoin
Set-TextValue "D15" "91.03"

# Row 16 - now ShibaInu (was Uniswap)
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.000008366"
$ws.Range("E16").Value = "  +1.47%  "

# Row 17 - now Uniswap (was ShibaInu)
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D17" "6.133"
$ws.Range("E17").Value = "  +2.30%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "29.301.51"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "240.42"
$ws.Range("E19").Value = "  -1.70%  "

# Row 20 - now Avalanche (was WrappedliquidstakedEther2.0)
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D20" "13.21"
$ws.Range("E20").Value = "  -0.47%  "

# Row 21 - now WrappedliquidstakedEther2.0 (was Avalanche)
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D21" "2.122.53"
$ws.Range("E21").Value = "  -0.93%  "

# Row 22 - Dai
Set-TextValue "D22" "1.000"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23 - Chainlink
Set-TextValue "D23" "7.762"
$ws.Range("E23").Value = "  -2.15%  "

# Row 24 - BinanceUSD
Set-TextValue "D24" "0.9999"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25 - Stellar
Set-TextValue "D25" "0.1593"
$ws.Range("E25").Value = "  -1.74%  "

# Row 26 - Monero
Set-TextValue "D26" "162.67"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.029"
$ws.Range("E27").Value = "  -0.18%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "18.51"
$ws.Range("E28").Value = "  -0.67%  "

# Row 29 - PancakeSwap
Set-TextValue "D29" "1.503"
$ws.Range("E29").Value = "  -0.35%  "

# Row 30 - Filecoin
Set-TextValue "D30" "4.408"
$ws.Range("E30").Value = "  -0.33%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "4.318"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32 - Toncoin
Set-TextValue "D32" "1.247"
$ws.Range("E32").Value = "  -4.01%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.05358"
$ws.Range("E33").Value = "  +2.35%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.938"
$ws.Range("E34").Value = "  +0.43%  "

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.176"
$ws.Range("E35").Value = "  -0.17%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.7492"
$ws.Range("E36").Value = "  -3.70%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.684"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38 - VeChain
Set-TextValue "D38" "0.01876"
$ws.Range("E38").Value = "  +0.41%  "

# Row 39 - Maker
Set-TextValue "D39" "1.241.62"
$ws.Range("E39").Value = "  +6.09%  "

# Row 40 - MXToken
Set-TextValue "D40" "2.732"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41 - FraxShare
Set-TextValue "D41" "6.516"
$ws.Range("E41").Value = "  +1.31%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.8916"
$ws.Range("E42").Value = "  -0.07%  "

# Row 43 - Quant
Set-TextValue "D43" "109.17"
$ws.Range("E43").Value = "  +4.45%  "

# Row 44 - Aave
Set-TextValue "D44" "72.36"
$ws.Range("E44").Value = "  -1.97%  "

# Row 45 - PaxDollar
Set-TextValue "D45" "0.9999"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46 - RocketPoolETH
Set-TextValue "D46" "2.024.31"
$ws.Range("E46").Value = "  -0.26%  "

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  +5.34%  "

# Row 48 - Mantle
Set-TextValue "D48" "0.5199"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49 - RenderToken
Set-TextValue "D49" "1.791"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  +0.15%  "

# Row 51 - TheSandbox
Set-TextValue "D51" "0.4338"
$ws.Range("E51").Value = "  +0.42%  "
